$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "1+1="
$ws.Range("I2").Value = "abc"
$ws.Range("B4").Value = "12 + 1 = `n2 3  4  333 `nx3 + 2x2 - 1"
$ws.Range("I4").Value = "'1"
